$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ARC")
$ws2 = $wb.Worksheets.Item("ARC_Volumes")

# ---- Sheet: ARC ----
$ws1.Range("G1").Value = @'
AttributeDescriptionEN
'@
$ws1.Range("H1").Value = @'
AttributeName
'@
$ws1.Range("G6").Value = @'
Recommendation: \n
- Room Identifier: Persistent Unique Numerical Designation \n
- Definition: A distinct, enduring numerical code assigned to each room for unambiguous identification within a building or complex. E.g., #007 \n
\n
Key Principles for Effective Data Management:
1. Persistence: Maintain the same identifier as long as the room's core attributes remain unchanged: \n
   - Position (not relocated to a different part of the building)\n
   - General size (minor adjustments are acceptable)\n
   - Primary function (e.g., remains an office or meeting room)
2. Non-reuse: Retire old identifiers rather than reassigning them to new spaces:
   - Prevents confusion in historical data analysis
   - Avoids potential errors in facility management systems
3. Consistency: Apply identifiers systematically across all building documentation and systems:
   - Floor plans
   - Building Information Models (BIM)
   - Facility management software
   - Maintenance records
4. Format Standardization: Adopt a consistent format for all identifiers:
   - E.g., three-digit numbers preceded by # (#001, #002, etc.)
   - Consider building or floor prefixes for larger complexes (B1-#007, B2-#007)
5. Change Management: Implement a formal process for identifier changes when necessary:
   - Document reasons for changes
   - Update all relevant systems and documentation
   - Maintain a change log for future reference
6. Resist Client-Driven Changes: If clients request room numbering changes during project execution:
   - Explain that the current identifiers are technical room numbers crucial for data integrity
   - Offer to add a separate, non-technical room number attribute to accommodate their preference
   - Emphasize the importance of maintaining the original identifiers for system consistency and future operations
Remember: While it may seem tempting to reorganize or reuse identifiers, maintaining consistency over time is crucial for reliable long-term data management and analysis in facility operations. Technical room numbers should remain stable, even if additional labeling schemes are introduced for client use.
'@
$ws1.Range("G7").Value = @'
Room Type Naming: A Key Element in Building Information Management
Definition: A type name for the room, e.g., Office, Meeting Room, Laboratory, Storage.
Best Practices:
1. Use clear, descriptive names (e.g., "Open Plan Office" instead of just "Office")
2. Establish a standardized list of room types for your organization
3. Avoid abbreviations or codes that may be unclear to some users
4. Consider including subtypes for more detailed classification (e.g., "Meeting Room - Large")
5. Regularly review and update naming conventions to ensure they meet evolving needs
Remember: Consistent and thoughtful room type naming is foundational for effective building information management. It supports improved decision-making, efficiency, and long-term data usability across various aspects of facility management, energy modeling, space utilization, and regulatory compliance.
'@
$ws1.Range("C8").Value = @'
Useable Space
'@
$ws1.Range("D8").Value = @'
Model each usable space as one volume per story without any overlaps:_x005F_x000D_
_x005F_x000D_
- Ideal: From structural slab top to next structural slab bottom_x005F_x000D_
_x005F_x000D_
- Alternative: From usable floor top to structural slab bottom_x005F_x000D_
_x005F_x000D_
- Horizontal extent: From interior wall to interior wall, to capture the entire usable volume_x005F_x000D_
_x005F_x000D_
- Fill every "Void" in the building. e.g. Shafts, Elevators, 
'@
$ws1.Range("G8").Value = @'
In the Attribute, you can define if it's an EXTERNAL or INTERNAL Space.
'@
$ws1.Range("H8").Value = @'
PredefinedType
'@
$ws1.Range("K8").Value = @'
EXTERNAL, INTERNAL
'@
$ws1.Range("F9").Value = @'
Pset_SpaceCommon
'@
$ws1.Range("G9").Value = @'
The attribute defines space location:
True: Outside (e.g., balcony)
False: Inside (e.g., kitchen, parking garage)
'@
$ws1.Range("H9").Value = @'
IsInteriorOrExteriorSpace
'@
$ws1.Range("G10").Value = @'
The attribute defines space location:
True: Outside (e.g., balcony)
False: Inside (e.g., kitchen, parking garage)
'@
$ws1.Range("H10").Value = @'
IsExternal
'@
$ws1.Range("K10").Value = @'
TRUE, FALSE
'@
$ws1.Range("C11").Value = @'
Void in double height rooms
'@
$ws1.Range("D11").Value = @'
Model "Void" above usable space:_x005F_x000D_
- One volume per story_x005F_x000D_
- Ensure contact with space below/above
'@
$ws1.Range("F11").ClearContents()
$ws1.Range("G11").Value = @'
Recommendation:
- Name them consistently for simple filtering. E.g. Void
'@
$ws1.Range("H11").Value = @'
LongName
'@
$ws1.Range("K11").Value = @'
test
'@
$ws1.Range("C12").Value = @'
Window
'@
$ws1.Range("E12").Value = @'
IfcWindow
'@
$ws1.Range("G12").Value = @'
Recommendation:
- Specify window type via PredefinedType
'@
$ws1.Range("K12").ClearContents()
$ws1.Range("G13").Value = @'
Recommendation:
- Specify door type via PredefinedType
'@
$ws1.Range("G14").Value = @'
Recommendation:
- Specify window type via PredefinedType
'@
$ws1.Range("K14").Value = @'
DOOR, GATE, TRAPDOOR, USERDEFINED, NOTDEFINED
'@
$ws1.Range("C15").Value = @'
Doors
'@
$ws1.Range("E15").Value = @'
IfcDoor
'@
$ws1.Range("G15").Value = @'
Recommendation:
- Specify door type via PredefinedType
'@

# ---- Sheet: ARC_Volumes ----
$ws2.Range("G1").Value = @'
AttributeDescriptionEN
'@
$ws2.Range("H1").Value = @'
AttributeName
'@

